$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new currency row (TRY / 949) right after the existing data
$ws.Range("A9").Value = "TRY"
$ws.Range("B9").Value = 949

# Resize the worksheet table (ListObject) to include the new row
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:C9"))

# Update the active selection to match the new state
$ws.Range("A8").Select()
